# Basic Bass model for fertilisers
# Adds the CLEAFS (CLean Energy Adoption in Fertiliser Sector) module to the
# FTT_variables workbook: a new "CLEAFS" worksheet (inserted right before
# "Time_Horizons") describing the four new variables (FERTD, MFERTD, BFTC,
# AQR), plus the corresponding time-horizon lookup rows on the
# "Time_Horizons" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "CLEAFS" worksheet immediately before "Time_Horizons"
# ---------------------------------------------------------------------
$timeHorizons = $wb.Worksheets.Item("Time_Horizons")
$cleafs = $wb.Worksheets.Add($timeHorizons)
$cleafs.Name = "CLEAFS"

# Header row (matches the header used on every other FTT variable sheet)
$cleafs.Range("A1").Value = "Variable name"
$cleafs.Range("B1").Value = "Read in?"
$cleafs.Range("C1").Value = "Code"
$cleafs.Range("D1").Value = "Description"
$cleafs.Range("E1").Value = "RowDim"
$cleafs.Range("F1").Value = "ColDim"
$cleafs.Range("G1").Value = "3DDim"
$cleafs.Range("H1").Value = "Conversion?"
$cleafs.Range("I1").Value = "Scenario"

# Row 2 - FERTD
$cleafs.Range("A2").Value = "FERTD"
$cleafs.Range("B2").Value = 1
$cleafs.Range("C2").Value = 3601000
$cleafs.Range("D2").Value = "CLEAFS Fertiliser demand (tN/year)"
$cleafs.Range("E2").Value = "TFTI"
$cleafs.Range("F2").Value = "TIME"
$cleafs.Range("G2").Value = "RSHORTTI"
$cleafs.Range("H2").Value = 0
$cleafs.Range("I2").Value = "S0"

# Row 3 - MFERTD
$cleafs.Range("A3").Value = "MFERTD"
$cleafs.Range("B3").Value = 1
$cleafs.Range("C3").Value = 3602000
$cleafs.Range("D3").Value = "CLEAFS Maximum potential fertiliser demand (tN/year)"
$cleafs.Range("E3").Value = "TFTI"
$cleafs.Range("F3").Value = "TIME"
$cleafs.Range("G3").Value = "RSHORTTI"
$cleafs.Range("H3").Value = 0
$cleafs.Range("I3").Value = "S0"

# Row 4 - BFTC
$cleafs.Range("A4").Value = "BFTC"
$cleafs.Range("B4").Value = 1
$cleafs.Range("C4").Value = 3603000
$cleafs.Range("D4").Value = "CLEAFS matrix of technology costs"
$cleafs.Range("E4").Value = "TFTI"
$cleafs.Range("F4").Value = "CFTI"
$cleafs.Range("G4").Value = "RSHORTTI"
$cleafs.Range("H4").Value = 0
$cleafs.Range("I4").Value = "S0"

# Row 5 - AQR
$cleafs.Range("A5").Value = "AQR"
$cleafs.Range("B5").Value = 1
$cleafs.Range("C5").Value = 3604000
$cleafs.Range("D5").Value = "Agriculture output projection"
$cleafs.Range("E5").Value = "RSHORTTI"
$cleafs.Range("F5").Value = "TIME"
$cleafs.Range("G5").Value = "NA"
$cleafs.Range("H5").Value = 0
$cleafs.Range("I5").Value = "S0"

# ---------------------------------------------------------------------
# 2. Add the matching lookup rows on "Time_Horizons"
# ---------------------------------------------------------------------
# NOTE: re-fetch the worksheet by name - after Worksheets.Add(Before:=...)
# the original $timeHorizons reference now resolves to whatever sheet
# occupies that former position (the newly inserted CLEAFS sheet), not
# the Time_Horizons sheet itself.
$timeHorizons = $wb.Worksheets.Item("Time_Horizons")
$timeHorizons.Range("A67").Value = "FERTD"
$timeHorizons.Range("B67").Value = "tl_1960"
$timeHorizons.Range("A68").Value = "MFERTD"
$timeHorizons.Range("B68").Value = "tl_1960"
$timeHorizons.Range("A69").Value = "AQR"
$timeHorizons.Range("B69").Value = "tl_2010"

# ---------------------------------------------------------------------
# 3. Restore/replicate cell-selection state on the other sheets
# ---------------------------------------------------------------------
$ftt_p = $wb.Worksheets.Item("FTT-P")
[void]$ftt_p.Activate()
$excel.ActiveWindow.ScrollRow = 1
[void]$ftt_p.Range("A1:I19").Select()

$ftt_tr = $wb.Worksheets.Item("FTT-Tr")
[void]$ftt_tr.Activate()
$excel.ActiveWindow.ScrollRow = 1
[void]$ftt_tr.Range("E20").Select()

$ftt_h = $wb.Worksheets.Item("FTT-H")
[void]$ftt_h.Activate()
[void]$ftt_h.Range("D2").Select()

# FTT-S is untouched by this commit - leave its view/selection as-is.

$ftt_fr = $wb.Worksheets.Item("FTT-Fr")
[void]$ftt_fr.Activate()
$excel.ActiveWindow.ScrollRow = 1
[void]$ftt_fr.Range("C17").Select()

[void]$timeHorizons.Activate()
$excel.ActiveWindow.ScrollRow = 40
[void]$timeHorizons.Range("E65").Select()

# CLEAFS is the sheet that should end up active/selected, matching the
# original author's final view state.
[void]$cleafs.Activate()
[void]$cleafs.Range("G6").Select()
